$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Fold_1")
$ws.Range("B2").Value = 5.45195
$ws.Range("C2").Value = 5.89945
$ws.Range("D2").Value = 37.50169999999999
$ws.Range("E2").Value = 0.17115
$ws.Range("F2").Value = 2219.401700000001
$ws.Range("G2").Value = 2014.1057
$ws.Range("H2").Value = 205.29605
$ws.Range("I2").Value = 205.29605
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 2798.105
$ws.Range("L2").Value = 2014.125
$ws.Range("M2").Value = 783.98
$ws.Range("N2").Value = 345.6012
$ws.Range("O2").Value = 438.3792
$ws.Range("B3").Value = 5.789
$ws.Range("C3").Value = 9.571999999999999
$ws.Range("D3").Value = 28.079
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 2327.595
$ws.Range("G3").Value = 2122.899
$ws.Range("H3").Value = 204.696
$ws.Range("I3").Value = 193.677
$ws.Range("J3").Value = 11.019
$ws.Range("K3").Value = 2476.9342
$ws.Range("L3").Value = 2122.95
$ws.Range("M3").Value = 353.9842
$ws.Range("N3").Value = 308.7674
$ws.Range("O3").Value = 45.2168
$ws.Range("B4").Value = 6.746
$ws.Range("C4").Value = 8.249000000000001
$ws.Range("D4").Value = 41.042
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 2510.203
$ws.Range("G4").Value = 2280.003
$ws.Range("H4").Value = 230.2
$ws.Range("I4").Value = 64.40900000000001
$ws.Range("J4").Value = 165.791
$ws.Range("K4").Value = 2467.796
$ws.Range("L4").Value = 2280.003
$ws.Range("M4").Value = 187.793
$ws.Range("N4").Value = 176.6258
$ws.Range("O4").Value = 11.1674

$ws = $wb.Worksheets.Item("Fold_2")
$ws.Range("B2").Value = 5.591900000000001
$ws.Range("C2").Value = 6.396350000000001
$ws.Range("D2").Value = 36.86314999999999
$ws.Range("E2").Value = 0.1292
$ws.Range("F2").Value = 2253.7439
$ws.Range("G2").Value = 2044.3985
$ws.Range("H2").Value = 209.34555
$ws.Range("I2").Value = 209.34555
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 2285.459
$ws.Range("L2").Value = 2044.411
$ws.Range("M2").Value = 241.048
$ws.Range("N2").Value = 241.048
$ws.Range("O2").Value = 0
$ws.Range("B3").Value = 5.963
$ws.Range("C3").Value = 9.461
$ws.Range("D3").Value = 28.427
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 2359.512
$ws.Range("G3").Value = 2141.396
$ws.Range("H3").Value = 218.116
$ws.Range("I3").Value = 201.182
$ws.Range("J3").Value = 16.934
$ws.Range("K3").Value = 2347.4208
$ws.Range("L3").Value = 2141.429
$ws.Range("M3").Value = 205.9918
$ws.Range("N3").Value = 205.9918
$ws.Range("O3").Value = 0
$ws.Range("B4").Value = 6.7
$ws.Range("C4").Value = 8.093999999999999
$ws.Range("D4").Value = 47.574
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 2580.234
$ws.Range("G4").Value = 2317.715
$ws.Range("H4").Value = 262.519
$ws.Range("I4").Value = 54.052
$ws.Range("J4").Value = 208.467
$ws.Range("K4").Value = 2395.058
$ws.Range("L4").Value = 2317.715
$ws.Range("M4").Value = 77.343
$ws.Range("N4").Value = 77.343
$ws.Range("O4").Value = 0

$ws = $wb.Worksheets.Item("Fold_3")
$ws.Range("B2").Value = 5.43975
$ws.Range("C2").Value = 6.10355
$ws.Range("D2").Value = 37.0707
$ws.Range("E2").Value = 0.17115
$ws.Range("F2").Value = 2222.8903
$ws.Range("G2").Value = 2017.33845
$ws.Range("H2").Value = 205.55185
$ws.Range("I2").Value = 205.55185
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 2804.228
$ws.Range("L2").Value = 2017.357
$ws.Range("M2").Value = 786.8710000000001
$ws.Range("N2").Value = 359.913
$ws.Range("O2").Value = 426.9582
$ws.Range("B3").Value = 5.483
$ws.Range("C3").Value = 7.719
$ws.Range("D3").Value = 32.87
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 2320.004
$ws.Range("G3").Value = 2051.97
$ws.Range("H3").Value = 268.033
$ws.Range("I3").Value = 220.571
$ws.Range("J3").Value = 47.462
$ws.Range("K3").Value = 2618.863
$ws.Range("L3").Value = 2051.924
$ws.Range("M3").Value = 566.939
$ws.Range("N3").Value = 353.6808
$ws.Range("O3").Value = 213.2584
$ws.Range("B4").Value = 6.231
$ws.Range("C4").Value = 8.112
$ws.Range("D4").Value = 43.556
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 2573.291
$ws.Range("G4").Value = 2233.933
$ws.Range("H4").Value = 339.358
$ws.Range("I4").Value = 119.435
$ws.Range("J4").Value = 219.923
$ws.Range("K4").Value = 2490.067
$ws.Range("L4").Value = 2233.933
$ws.Range("M4").Value = 256.134
$ws.Range("N4").Value = 190.343
$ws.Range("O4").Value = 65.791

$ws = $wb.Worksheets.Item("Fold_4")
$ws.Range("B2").Value = 5.734050000000002
$ws.Range("C2").Value = 6.146849999999999
$ws.Range("D2").Value = 39.35619999999999
$ws.Range("E2").Value = 0.04355
$ws.Range("F2").Value = 2295.672100000001
$ws.Range("G2").Value = 2069.291249999999
$ws.Range("H2").Value = 226.3809
$ws.Range("I2").Value = 226.3809
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 2185.4816
$ws.Range("L2").Value = 2069.311
$ws.Range("M2").Value = 116.1706
$ws.Range("N2").Value = 80.70360000000001
$ws.Range("O2").Value = 35.467
$ws.Range("B3").Value = 6.113
$ws.Range("C3").Value = 9.366
$ws.Range("D3").Value = 29.606
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 2390.136
$ws.Range("G3").Value = 2164.023
$ws.Range("H3").Value = 226.112
$ws.Range("I3").Value = 213.845
$ws.Range("J3").Value = 12.267
$ws.Range("K3").Value = 2231.9788
$ws.Range("L3").Value = 2163.984
$ws.Range("M3").Value = 67.9948
$ws.Range("N3").Value = 67.9948
$ws.Range("O3").Value = 0
$ws.Range("B4").Value = 6.7
$ws.Range("C4").Value = 8.093999999999999
$ws.Range("D4").Value = 47.574
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 2580.234
$ws.Range("G4").Value = 2317.715
$ws.Range("H4").Value = 262.519
$ws.Range("I4").Value = 54.052
$ws.Range("J4").Value = 208.467
$ws.Range("K4").Value = 2328.9178
$ws.Range("L4").Value = 2317.715
$ws.Range("M4").Value = 11.2028
$ws.Range("N4").Value = 11.2028
$ws.Range("O4").Value = 0

$ws = $wb.Worksheets.Item("Fold_5")
$ws.Range("B2").Value = 5.656750000000001
$ws.Range("C2").Value = 5.934599999999999
$ws.Range("D2").Value = 39.77945
$ws.Range("E2").Value = 0.16955
$ws.Range("F2").Value = 2266.771
$ws.Range("G2").Value = 2056.2701
$ws.Range("H2").Value = 210.50085
$ws.Range("I2").Value = 210.50085
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 2393.0558
$ws.Range("L2").Value = 2056.285
$ws.Range("M2").Value = 336.7708
$ws.Range("N2").Value = 152.2866
$ws.Range("O2").Value = 184.4844
$ws.Range("B3").Value = 6.099
$ws.Range("C3").Value = 7.875
$ws.Range("D3").Value = 33.748
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 2371.432
$ws.Range("G3").Value = 2135.853
$ws.Range("H3").Value = 235.579
$ws.Range("I3").Value = 199.385
$ws.Range("J3").Value = 36.195
$ws.Range("K3").Value = 2332.4364
$ws.Range("L3").Value = 2135.884
$ws.Range("M3").Value = 196.5524
$ws.Range("N3").Value = 126.9988
$ws.Range("O3").Value = 69.5534
$ws.Range("B4").Value = 7.162
$ws.Range("C4").Value = 7.283
$ws.Range("D4").Value = 49.794
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 2576.937
$ws.Range("G4").Value = 2356.575
$ws.Range("H4").Value = 220.362
$ws.Range("I4").Value = 121.688
$ws.Range("J4").Value = 98.675
$ws.Range("K4").Value = 2467.0234
$ws.Range("L4").Value = 2356.575
$ws.Range("M4").Value = 110.4484
$ws.Range("N4").Value = 44.6264
$ws.Range("O4").Value = 65.8222
